# Append two new translation key/value rows to the "translations" sheet,
# fixing #2006, #2096, #2106.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

# Find the first empty row right after the existing data (row 68 -> 69).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$ws.Cells.Item($newRow1, 1).Value2 = "contribution_type_journal"
$ws.Cells.Item($newRow1, 2).Value2 = "Zeitschriftenbeitrag"

$ws.Cells.Item($newRow2, 1).Value2 = "contribution_type_book"
$ws.Cells.Item($newRow2, 2).Value2 = "Buchbeitrag"

Write-Host "Added rows $newRow1 and $newRow2 to sheet '$($ws.Name)'"
